# Updates cryptos list values (price/volume columns) to match the
# Sat Jul 22 11:08:06 UTC 2023 GitHub Actions refresh, including the
# Polygon / WrappedEther row swap (rows 12-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = '''29.914.69'
$ws.Cells.Item(2,5).Value = '  +0.36%  '
$ws.Cells.Item(3,4).Value = '''1.889.17'
$ws.Cells.Item(3,5).Value = '  +0.02%  '
$ws.Cells.Item(4,5).Value = '  -0.03%  '
$ws.Cells.Item(5,4).Value = '''0.7737'
$ws.Cells.Item(5,5).Value = '  +0.59%  '
$ws.Cells.Item(6,4).Value = '''243.00'
$ws.Cells.Item(6,5).Value = '  -0.52%  '
$ws.Cells.Item(7,4).Value = '''1.000'
$ws.Cells.Item(7,5).Value = '  -0.03%  '
$ws.Cells.Item(8,4).Value = '''0.3107'
$ws.Cells.Item(8,5).Value = '  -0.39%  '
$ws.Cells.Item(9,4).Value = '''25.67'
$ws.Cells.Item(9,5).Value = '  +1.84%  '
$ws.Cells.Item(10,4).Value = '''0.07168'
$ws.Cells.Item(10,5).Value = '  -0.66%  '
$ws.Cells.Item(11,4).Value = '''0.08577'
$ws.Cells.Item(11,5).Value = '  +5.94%  '
$ws.Cells.Item(12,2).Value = 'WrappedEther'
$ws.Cells.Item(12,3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(12,4).Value = '''1.997.03'
$ws.Cells.Item(12,5).Value = '  +3.05%  '
$ws.Cells.Item(13,2).Value = 'Polygon'
$ws.Cells.Item(13,3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(13,4).Value = '''0.7642'
$ws.Cells.Item(13,5).Value = '  +0.05%  '
$ws.Cells.Item(14,4).Value = '''5.365'
$ws.Cells.Item(14,5).Value = '  -2.34%  '
$ws.Cells.Item(15,4).Value = '''93.87'
$ws.Cells.Item(15,5).Value = '  +1.76%  '
$ws.Cells.Item(16,4).Value = '''6.188'
$ws.Cells.Item(16,5).Value = '  +0.78%  '
$ws.Cells.Item(17,4).Value = '''30.040.26'
$ws.Cells.Item(17,5).Value = '  +0.75%  '
$ws.Cells.Item(18,4).Value = '''13.77'
$ws.Cells.Item(18,5).Value = '  -0.81%  '
$ws.Cells.Item(19,4).Value = '''244.58'
$ws.Cells.Item(19,5).Value = '  +0.78%  '
$ws.Cells.Item(20,4).Value = '''0.000007816'
$ws.Cells.Item(20,5).Value = '  +0.67%  '
$ws.Cells.Item(21,4).Value = '''2.213.80'
$ws.Cells.Item(21,5).Value = '  +2.20%  '
$ws.Cells.Item(22,4).Value = '''0.9972'
$ws.Cells.Item(23,4).Value = '''7.960'
$ws.Cells.Item(23,5).Value = '  -2.90%  '
$ws.Cells.Item(25,5).Value = '  +4.88%  '
$ws.Cells.Item(26,4).Value = '''9.358'
$ws.Cells.Item(26,5).Value = '  -0.78%  '
$ws.Cells.Item(27,4).Value = '''162.02'
$ws.Cells.Item(27,5).Value = '  -0.17%  '
$ws.Cells.Item(28,4).Value = '''18.75'
$ws.Cells.Item(28,5).Value = '  -0.02%  '
$ws.Cells.Item(29,5).Value = '  -0.23%  '
$ws.Cells.Item(30,4).Value = '''1.442'
$ws.Cells.Item(30,5).Value = '  -1.42%  '
$ws.Cells.Item(31,4).Value = '''1.541'
$ws.Cells.Item(31,5).Value = '  -0.57%  '
$ws.Cells.Item(32,4).Value = '''4.510'
$ws.Cells.Item(32,5).Value = '  +1.41%  '
$ws.Cells.Item(33,4).Value = '''4.112'
$ws.Cells.Item(33,5).Value = '  +0.91%  '
$ws.Cells.Item(34,4).Value = '''0.05435'
$ws.Cells.Item(34,5).Value = '  -1.53%  '
$ws.Cells.Item(35,4).Value = '''1.241'
$ws.Cells.Item(35,5).Value = '  -1.11%  '
$ws.Cells.Item(36,4).Value = '''0.7470'
$ws.Cells.Item(36,5).Value = '  -0.05%  '
$ws.Cells.Item(37,5).Value = '  +0.06%  '
$ws.Cells.Item(38,4).Value = '''2.695'
$ws.Cells.Item(38,5).Value = '  +2.26%  '
$ws.Cells.Item(39,4).Value = '''0.01966'
$ws.Cells.Item(39,5).Value = '  +2.48%  '
$ws.Cells.Item(40,5).Value = '  +0.01%  '
$ws.Cells.Item(41,4).Value = '''0.4464'
$ws.Cells.Item(41,5).Value = '  +1.15%  '
$ws.Cells.Item(42,4).Value = '''1.106.97'
$ws.Cells.Item(42,5).Value = '  -3.22%  '
$ws.Cells.Item(43,4).Value = '''73.45'
$ws.Cells.Item(43,5).Value = '  -0.16%  '
$ws.Cells.Item(44,4).Value = '''6.095'
$ws.Cells.Item(44,5).Value = '  +3.23%  '
$ws.Cells.Item(45,4).Value = '''0.8504'
$ws.Cells.Item(45,5).Value = '  +0.40%  '
$ws.Cells.Item(46,5).Value = '  -0.05%  '
$ws.Cells.Item(47,4).Value = '''103.44'
$ws.Cells.Item(47,5).Value = '  +0.57%  '
$ws.Cells.Item(48,4).Value = '''1.871'
$ws.Cells.Item(48,5).Value = '  -0.56%  '
$ws.Cells.Item(49,4).Value = '''7.618'
$ws.Cells.Item(49,5).Value = '  +2.53%  '
$ws.Cells.Item(50,4).Value = '''2.078.84'
$ws.Cells.Item(50,5).Value = '  +1.19%  '
$ws.Cells.Item(51,4).Value = '''2.994'
$ws.Cells.Item(51,5).Value = '  -0.66%  '
